$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9245747923851013
$ws.Range("B1").Value = 1.55974280834198
$ws.Range("D1").Value = 1.61381733417511
$ws.Range("E1").Value = 1.054721474647522
